# Samuel's branch change:
#  - Move the selection on Sheet1 from B2 to B9 (Sheet1 stops being the active/tab-selected sheet).
#  - Add a new "Sheet2" right after "Sheet1", make it the active sheet/tab.
#  - Put the text "Cambio Branch de Samuel" in Sheet2!A1.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Update the selection on Sheet1 (was B2, now B9) before we hand focus to the new sheet.
[void]$sheet1.Range("B9").Select()

# Insert the new sheet immediately after Sheet1 (not using the default Add(),
# which would insert before the active sheet).
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Sheet2"

# Content for the new sheet.
$newSheet.Range("A1").Value = "Cambio Branch de Samuel"

# Match the original page margins used elsewhere in the workbook
# (PageSetup margins are expressed in points, 72 points per inch).
$newSheet.PageSetup.LeftMargin = 0.75 * 72
$newSheet.PageSetup.RightMargin = 0.75 * 72
$newSheet.PageSetup.TopMargin = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

# Sheet2 becomes the active/selected tab (activeTab=1, tabSelected on Sheet2).
[void]$newSheet.Activate()
